$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data row by row.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.578.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.01%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.452.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.33%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.46%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.454.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.42%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.95%  "

# Row 10
$ws.Range("E10").Value = "  +0.37%  "

# Row 11
$ws.Range("E11").Value = "  +3.37%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.389"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.64%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.042.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.36%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.27%  "

# Row 15
$ws.Range("E15").Value = "  -0.85%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.95%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.450.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.42%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.719.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "

# Row 19
$ws.Range("E19").Value = "  +9.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.86%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.40%  "

# Row 23
$ws.Range("E23").Value = "  +3.54%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.48%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.25%  "

# Row 27
$ws.Range("E27").Value = "  +0.57%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.590.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.180"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.50%  "

# Row 31
$ws.Range("E31").Value = "  +0.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.14%  "

# Row 33
$ws.Range("E33").Value = "  +2.61%  "

# Row 34
$ws.Range("E34").Value = "  -11.48%  "

# Row 35
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("E36").Value = "  +3.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.479.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.50%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.50%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.23%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.01%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "166.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.42%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "28.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.05%  "

# Row 43
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0783"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.54%  "

# Row 44
$ws.Range("E44").Value = "  +3.91%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.89%  "

# Row 46
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("E47").Value = "  +4.32%  "

# Row 48
$ws.Range("E48").Value = "  +3.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.583.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.96%  "

# Row 50
$ws.Range("E50").Value = "  -1.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.72%  "
